# ---------------------------------------------------------------------------
# kra_master_database.xlsx edit
#   - camelCase field rollout: add Merged_From_Count, Merge_Sources,
#     Best_Score, Pre-Amount, date, pin, taxpayerName, preAmount,
#     finalAmount, year, officerName, station (columns J:U) on
#     KRA_Database.
#   - Re-home the dedup/merge record (James Mutoro Kitui) to row 2 with its
#     merge metadata, shift the two "Test" rows down to 3/4.
#   - Append two new rows: a pre-amount extraction test row (12) and a
#     camelCase-field test row (13).
#   - Normalise D11 (Year) to a real number.
#   - Refresh the Database_Summary rollup numbers.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("KRA_Database")
$ws2 = $wb.Worksheets.Item("Database_Summary")

# --- helper: write a string value while forcing TEXT storage, so Excel's
#     auto-detection doesn't silently turn date-/number-looking strings
#     (e.g. "2024-01-15", "14,769.50") into numeric/date serials.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ===========================================================================
# 1. New header row (J1:U1) — camelCase + merge-metadata columns
# ===========================================================================
$newHeaders = @(
    "Merged_From_Count", "Merge_Sources", "Best_Score", "Pre-Amount",
    "date", "pin", "taxpayerName", "preAmount", "finalAmount",
    "year", "officerName", "station"
)
$newHeaderCols = @("J","K","L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Range($newHeaderCols[$i] + "1").Value = $newHeaders[$i]
}

# Match the existing header styling (bold, centered, thin border) by
# copying the format from an existing header cell onto the new ones.
$ws.Range("A1").Copy()
$ws.Range("J1:U1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ===========================================================================
# 2. Row 2 — the deduplicated / merged record (previously row 4's data:
#    James Mutoro Kitui), now carrying merge metadata in J:L
# ===========================================================================
$ws.Range("A2").Value = "04th September, 2025"
$ws.Range("B2").Value = "A005977112Z"
$ws.Range("C2").Value = "James Mutoro Kitui"
$ws.Range("D2").Value = 2024
$ws.Range("E2").Value = "Franciscar Nyangweta"
$ws.Range("F2").Value = "KITALE"
$ws.Range("G2").Value = "2025-09-21 22:53:29"
$ws.Range("H2").Value = "multi_format_extractor"
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = "Unknown, Unknown"
$ws.Range("L2").Value = 100

# ===========================================================================
# 3. Row 3 — the "Test User" sample record (previously row 2's data)
# ===========================================================================
Set-TextValue $ws.Range("A3") "2024-01-15"
$ws.Range("B3").Value = "A123456789X"
$ws.Range("C3").Value = "Test User"
$ws.Range("D3").Value = 2023
$ws.Range("E3").Value = "Test Officer"
$ws.Range("F3").Value = "Test Station"
$ws.Range("G3").Value = "2025-09-21 22:51:45"
$ws.Range("H3").Value = "test_app"
$ws.Range("I3").Value = 1

# ===========================================================================
# 4. Row 4 — Ezekiel Kipserem Korir (previously row 3's data)
# ===========================================================================
$ws.Range("A4").Value = "4th September, 2025"
$ws.Range("B4").Value = "A009775891W"
$ws.Range("C4").Value = "Ezekiel Kipserem Korir"
$ws.Range("D4").Value = 2024
$ws.Range("E4").Value = "Franciscar Nyangweta"
$ws.Range("F4").Value = "KITALE"
$ws.Range("G4").Value = "2025-09-21 22:53:29"
$ws.Range("H4").Value = "multi_format_extractor"
$ws.Range("I4").Value = 2

# ===========================================================================
# 5. Rows 5-10 are unchanged — only row 11 needs its Year cell converted
#    from text to a genuine number.
# ===========================================================================
$ws.Range("D11").Value = 2024

# ===========================================================================
# 6. Row 12 (NEW) — pre-amount extractor test record
# ===========================================================================
Set-TextValue $ws.Range("A12") "2024-09-15"
$ws.Range("B12").Value = "A123456789X"
$ws.Range("C12").Value = "Peter Kimutai Telengech"
$ws.Range("D12").Value = 2024
$ws.Range("E12").Value = "John Doe"
$ws.Range("F12").Value = "NAIROBI"
$ws.Range("G12").Value = "2025-09-22 08:40:33"
$ws.Range("H12").Value = "test_pre_amount"
$ws.Range("I12").Value = 11
Set-TextValue $ws.Range("M12") "14,769.50"

# ===========================================================================
# 7. Row 13 (NEW) — camelCase field test record (A:F left blank)
# ===========================================================================
$ws.Range("G13").Value = "2025-09-22 08:52:05"
$ws.Range("H13").Value = "camelCase_test"
$ws.Range("I13").Value = 12
Set-TextValue $ws.Range("N13") "2024-09-22"
$ws.Range("O13").Value = "A123456789X"
$ws.Range("P13").Value = "John Doe Test"
Set-TextValue $ws.Range("Q13") "14,769.50"
# R13 (finalAmount) intentionally left blank, as requested
Set-TextValue $ws.Range("S13") "2024"
$ws.Range("T13").Value = "Test Officer"
$ws.Range("U13").Value = "NAIROBI"

# ===========================================================================
# 8. Database_Summary rollup refresh
# ===========================================================================
$ws2.Range("B2").Value = 12
$ws2.Range("B3").Value = "2025-09-22 08:52:05"
$ws2.Range("B6").ClearContents()
Set-TextValue $ws2.Range("B7") "2024-09-22"
$ws2.Range("B8").Value = 2
$ws2.Range("B9").Value = 2
